# VAC_CRIANCA.xlsx edit:
#  - Introduce a new prompt type "adate" (replacing "custom_date" usage
#    throughout the survey + model sheets) with comment text describing
#    the new date-save behaviour.
#  - Add the new "adate" row to the prompt_types sheet.
#  - Update selection / active-tab bookkeeping: the "model" sheet becomes
#    the active tab (was "survey").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) survey sheet: every "type" cell (column D) that referenced
#    "custom_date" now references "adate".
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$surveyRows = @(9,14,20,25,47,52,58,63,71,76,82,87,110,115,121,126,134,139,145,150,171,176,182,187,195,200,206,211,228,233,239,244)
foreach ($r in $surveyRows) {
    $survey.Range("D$r").Value = "adate"
}

# ---------------------------------------------------------------------
# 2) model sheet: every "type" cell (column B) that referenced
#    "custom_date" now references "adate".
# ---------------------------------------------------------------------
$model = $wb.Worksheets.Item("model")
$modelRows = @(42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,163,167,171,175,179,183,187,191,195,199,203,207,211,215,219,227,233)
foreach ($r in $modelRows) {
    $model.Range("B$r").Value = "adate"
}

# ---------------------------------------------------------------------
# 3) prompt_types sheet: append the new "adate" definition row.
# ---------------------------------------------------------------------
$promptTypes = $wb.Worksheets.Item("prompt_types")
$promptTypes.Range("A4").Value = "adate"
$promptTypes.Range("B4").Value = "string"
$promptTypes.Range("C4").Value = "string"
$promptTypes.Range("D4").Value = "Save only mm.dd.yyyy with support for ?? at all positions"

# ---------------------------------------------------------------------
# 4) View/selection bookkeeping to match the authored workbook state.
# ---------------------------------------------------------------------

# prompt_types: selection moves to C20.
$promptTypes.Activate()
$promptTypes.Range("C20").Select()

# survey: selection moves to D5, no longer the active tab.
$survey.Activate()
$survey.Range("D5").Select()

# model: selection moves to B45 and it becomes the active tab.
$model.Activate()
$model.Range("B45").Select()

Write-Host "edit applied"
